# "Update automatico" refresh: the timestamp log in column D shifts down
# one slot and a new timestamp (~2021-02-22 13:31:34) is stamped onto the
# most recent block of rows (2-15), pushing the two older timestamp
# blocks (rows 16-29, rows 30-43) to the values that used to belong to
# the block above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$blocks = @(
    @{ FirstRow = 2;  LastRow = 15; Value = 44249.56359322634 },
    @{ FirstRow = 16; LastRow = 29; Value = 44249.54226155092 },
    @{ FirstRow = 30; LastRow = 43; Value = 44249.52093606482 }
)

foreach ($block in $blocks) {
    for ($row = $block.FirstRow; $row -le $block.LastRow; $row++) {
        $ws.Cells.Item($row, 4).Value2 = $block.Value
    }
}
